$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H6").Value = "5c3046fa-bb76-4d07-b4d7-2cdfd85193a2"
$ws.Range("H7").Value = "4882c6b1-77f1-4cba-ba86-c26e062dab71"
$ws.Range("H8").Value = "e6139de2-edd1-4516-a7c0-8828e6a6d4b3"
$ws.Range("H9").Value = "948d5ddb-caef-4683-af04-a9936974d326"
$ws.Range("H10").Value = "7f0fd156-d853-438e-a643-62156d1d8f90"
$ws.Range("H11").Value = "81ac8b91-89c0-4ec2-89ec-fa5d8ee36bc3"
$ws.Range("H12").Value = "13a43ce9-da28-48a9-90aa-7a1ff58d7bab"
$ws.Range("H13").Value = "a638ff45-26b4-46d9-81b2-ce35f105277f"
$ws.Range("H14").Value = "9c1b5ed8-c102-41c4-a8a7-1e9a9b5168a4"
$ws.Range("H15").Value = "f96e38c5-a3e0-483e-ba36-5d5934b9503f"
$ws.Range("H16").Value = "74b39964-6432-456d-8471-338dcd6d42e3"
$ws.Range("H17").Value = "4000bd76-bd24-49e2-849a-839dabc51be8"
$ws.Range("H19").Value = "5cd2a955-9d69-42e9-aa54-03d3f2e1e1be"
$ws.Range("H20").Value = "e10797f7-b3cf-42c5-9ffc-f8fa529c89df"
$ws.Range("H21").Value = "35190394-6745-4cf5-ac25-7637e9c1d08d"
$ws.Range("H22").Value = "9e6b7969-254a-4e5a-a2b4-7844fb3bb0ec"
$ws.Range("H23").Value = "9a969bf1-b846-4f90-b8bc-995c5cbe0ca1"
$ws.Range("H24").Value = "20074e00-8db4-4ced-950d-28dfc6a384a6"
$ws.Range("H25").Value = "7c82d85a-a562-4a4c-be1d-eb3644f85357"
$ws.Range("H26").Value = "1f2f3f60-99f7-489c-a320-87ec93cd091e"
$ws.Range("H27").Value = "426b3fe4-013b-4d2e-a808-412399b989f4"
$ws.Range("H28").Value = "263b09ed-942a-4ab5-a387-7b65798bb7e1"
$ws.Range("H29").Value = "c52dbac0-2f05-461a-9f17-7460c1846953"
$ws.Range("H30").Value = "028f272c-4fb9-44b2-a253-ed1957e69daf"
$ws.Range("H31").Value = "7aecbcab-353a-4168-8661-21d361ca2c12"
$ws.Range("H32").Value = "1f05683a-c313-4d6c-ba3c-dab029472ba9"
$ws.Range("H33").Value = "a736fbae-3a61-44f9-9a4c-ed4720b91762"
$ws.Range("H34").Value = "76da67c9-4093-44cf-b3c4-5aa5dbe243f1"
$ws.Range("H35").Value = "bc1fc326-08ba-4d2b-a7c5-3115ffca70a0"
$ws.Range("H36").Value = "6acfa2f1-2f37-4bd9-bd69-166a0bf82363"
$ws.Range("H37").Value = "57d1b10a-da42-45ab-96c7-590470dc3f3e"
$ws.Range("H38").Value = "8b063687-59c5-4d4e-ab4b-6313e6c3f08f"
$ws.Range("H39").Value = "e2c6ea89-c5e0-4d20-96d1-a2f5cb64efa2"
$ws.Range("H40").Value = "3b3b86c0-dbbc-43ad-a21b-6f0704a8f8cc"
$ws.Range("H41").Value = "72b97a24-1c89-464d-a09a-71475050f108"
$ws.Range("H42").Value = "f9466b39-f04e-4376-897e-cac3889791a7"
$ws.Range("H43").Value = "15fb7579-253d-4914-9ccc-b506da18b409"
$ws.Range("H44").Value = "c4149881-51a1-4dc9-acd0-3b7908637b04"
$ws.Range("H45").Value = "0ac95d85-02ac-4bac-b097-dd5a41aa2936"
$ws.Range("H46").Value = "6b1b8963-6c34-4a38-9742-940c8e52c170"
$ws.Range("H47").Value = "df232722-5f3d-4282-ba45-e898f6d50fd7"
$ws.Range("H48").Value = "c703b4fc-96a9-4bb1-8501-d81233b0e695"
$ws.Range("H49").Value = "dfda3db0-757f-4520-b0fd-5c10769827e4"
$ws.Range("H50").Value = "1dd59022-ffd2-4940-8604-fc004d3d8f34"
$ws.Range("H51").Value = "eb72acdf-79a3-4f23-aa27-66b37674d5b6"
$ws.Range("H52").Value = "a0ba7912-9cce-41d2-8e9e-5d6f05b31c5c"
$ws.Range("H53").Value = "7f4ecbcd-d6b1-4ff3-80c5-333599aa0d36"
$ws.Range("H54").Value = "a205f55e-562b-4eb2-bc42-9faaf3d6a3ab"
$ws.Range("H55").Value = "7b56f30f-635a-41f1-91d6-522b2e939fdd"
$ws.Range("H56").Value = "72d661e9-d76d-468e-8347-ed598c30b2a3"
$ws.Range("H57").Value = "af2b1808-e6dd-49c5-9269-75664b0a3f75"
$ws.Range("H58").Value = "a163e47d-eb17-4caf-82f3-3eae56538625"
$ws.Range("H59").Value = "08e9bca4-b4d3-45fa-a5e2-d3e22f11c8d3"
$ws.Range("H60").Value = "ac914fee-17eb-46bb-8ff3-10a0c2922076"
$ws.Range("H61").Value = "9b7a113a-d186-42a6-a11b-38a2795dd6d4"
$ws.Range("H62").Value = "c9eb4f53-a13c-475e-92b2-cf6049001957"
$ws.Range("H63").Value = "e19f9abd-e119-4ff7-96bb-7d054beb18e0"
$ws.Range("H64").Value = "e384b13a-654a-475f-a362-f625dbf5e10c"
$ws.Range("H65").Value = "32372d09-9dd6-451c-b552-523c33d815ba"
$ws.Range("H66").Value = "2d1121a6-2e34-4075-b1ff-c7bc514924e2"
